$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$tr.Text = "Creating an R Project and using here() function`rData processing using DPLYR`rPlotting using GGPLOT`rLinear model (Regression) and creating dummy variables`rSome OH queries`r"
